# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap country names for a handful of rows (A column) ---
# These swaps reuse existing shared-string text so that the net effect
# matches the source diff (rows keep their statistics, only the
# displayed country name changes because the underlying shared string
# table entries were re-ordered upstream).
$ws.Range("A47").Value  = "Panama"
$ws.Range("A48").Value  = "Rumania"

$ws.Range("A202").Value = "Dominica"
$ws.Range("A203").Value = "Fiyi"

$ws.Range("A206").Value = "Groenlandia"
$ws.Range("A207").Value = "Islas Malvinas"

$ws.Range("A208").Value = "Islas Turcas y Caicos"
$ws.Range("A209").Value = "Santa Sede"

# --- Update the "last refreshed" timestamp label (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Junio de 2020 a las 02:54"

# --- Update statistics figures ---
# Row 4 - Estados Unidos
$ws.Range("B4").Value = 2263634
$ws.Range("C4").Value = 27907
$ws.Range("D4").Value = 930783
$ws.Range("E4").Value = 1212163
$ws.Range("G4").Value = 747
$ws.Range("H4").Value = 120688

# Row 37 - Argentina
$ws.Range("B37").Value = 37510
$ws.Range("C37").Value = 1958
$ws.Range("E37").Value = 25841
$ws.Range("G37").Value = 35
$ws.Range("H37").Value = 948

# Row 47 - now Panama
$ws.Range("B47").Value = 23351
$ws.Range("C47").Value = 754
$ws.Range("D47").Value = 13782
$ws.Range("E47").Value = 9094
$ws.Range("G47").Value = 5
$ws.Range("H47").Value = 475

# Row 48 - now Rumania
$ws.Range("B48").Value = 23080
$ws.Range("C48").Value = 320
$ws.Range("D48").Value = 16308
$ws.Range("E48").Value = 5299
$ws.Range("G48").Value = 22
$ws.Range("H48").Value = 1473

# Row 91 - Venezuela
$ws.Range("B91").Value = 3483
$ws.Range("C91").Value = 97
$ws.Range("E91").Value = 2620

# Row 191 - Gambia
$ws.Range("B191").Value = 36
$ws.Range("C191").Value = 2
$ws.Range("E191").Value = 11

# Row 197 - Granada
$ws.Range("D197").Value = 23
$ws.Range("E197").Value = 0

# Row 208 - now Islas Turcas y Caicos
$ws.Range("D208").Value = 11
$ws.Range("H208").Value = 1

# Row 209 - now Santa Sede
$ws.Range("D209").Value = 12
$ws.Range("H209").Value = 0
